$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.573.45'
$ws.Range('E2').Value = '  +2.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.606.68'
$ws.Range('E3').Value = '  +2.61%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.60'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('E6').Value = '  +4.36%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.87'
$ws.Range('E8').Value = '  +7.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.47'
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('E10').Value = '  +2.24%  '
$ws.Range('E11').Value = '  +2.49%  '
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.836.07'
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.620.07'
$ws.Range('E14').Value = '  +3.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.575.86'
$ws.Range('E15').Value = '  +2.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.537'
$ws.Range('E17').Value = '  +2.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.45'
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.58'
$ws.Range('E19').Value = '  +4.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.60'
$ws.Range('E20').Value = '  +3.40%  '
$ws.Range('E21').Value = '  +1.62%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  +1.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.22'
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.12'
$ws.Range('E26').Value = '  +1.81%  '
$ws.Range('E27').Value = '  +3.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.26'
$ws.Range('E28').Value = '  +3.08%  '
$ws.Range('E29').Value = '  +2.44%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  +2.41%  '
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('E34').Value = '  +3.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.408.89'
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('E37').Value = '  +4.56%  '
$ws.Range('E38').Value = '  +4.60%  '
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('E40').Value = '  +2.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.538'
$ws.Range('E41').Value = '  +3.72%  '
$ws.Range('E42').Value = '  +1.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0493'
$ws.Range('E43').Value = '  +7.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '53.53'
$ws.Range('E44').Value = '  +25.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.798'
$ws.Range('E45').Value = '  +3.37%  '
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.90'
$ws.Range('E47').Value = '  +3.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.29'
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.746.59'
$ws.Range('E49').Value = '  +2.40%  '
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '86.66'
$ws.Range('E51').Value = '  +1.84%  '
